# daily test level 조정
# Adjust daily-test level/score values, row heights, and the default
# (unused/beyond-data) column alignment for the UserList sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (level/score adjustments) ---
# user1 (row 5): level 1 -> 11, last test date -> 45438 (2024-05-01)
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 45438
# user2 (row 6): level 700 -> 710
$ws.Range("E6").Value = 710

# --- Row height adjustments: data rows 2-14 go from 21 to 19.5 ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}

# --- Default column alignment: general -> left for columns A-D and F ---
# (These columns carry a column-level default format - xf entries used
# beyond the populated data range - that switches from "general" to
# "left" horizontal alignment.)
$ws.Columns("A:D").HorizontalAlignment = -4131
$ws.Columns("F").HorizontalAlignment = -4131
